$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Text edit inside the "INSTALAÇÃO" paragraph: remove "isso " and add
#    a comma after "nuvem":
#    " ... na nuvem logo após isso deverá ..." ->
#    " ... na nuvem, logo após deverá ..."
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "na nuvem logo após isso deverá",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "na nuvem, logo após deverá", 2)

# ---------------------------------------------------------------------
# 2) Move the hidden "_GoBack" bookmark from the end of the paragraph
#    that was just edited to the end of the first paragraph (right
#    after "... recém-nascidos " and before the paragraph mark).
#
#    The emulator's Bookmarks.Add() snaps any collapsed range sitting
#    exactly on a paragraph boundary back to paragraph 1, so we work
#    around it by temporarily inserting a spacer character after the
#    target point (which moves the boundary elsewhere), anchoring the
#    bookmark just before the spacer, and then deleting the spacer --
#    the now-anchored bookmark stays put.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$goBackPos = $p1.Range.End - 1

$spacer = $d.Range($goBackPos, $goBackPos)
$spacer.InsertAfter("X")

$bmRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($goBackPos, $goBackPos + 1).Delete()
